$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data (D:K) to (E:L)
$ws.Columns("D").Insert()

# Copy formatting from the (now shifted) column E onto the new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 330900
$ws.Range("D9").Value = 121500
$ws.Range("D10").Value = 209400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 107100
$ws.Range("D17").Value = 251300
$ws.Range("D18").Value = 79500
$ws.Range("D20").Value = -100
$ws.Range("D21").Value = 186500
$ws.Range("D22").Value = 52200
$ws.Range("D23").Value = 27200
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 27200
$ws.Range("D27").Value = 19700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 100
$ws.Range("D33").Value = 19700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 19700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 48000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 9300
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 6700
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 39800
$ws.Range("D48").Value = 2039900
$ws.Range("D49").Value = 9800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 9300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2198300
$ws.Range("D57").Value = 46600
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 10300
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 1290800
$ws.Range("D62").Value = 100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1395300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -128800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 803000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 19700
$ws.Range("D83").Value = 107100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 136500
$ws.Range("D91").Value = -54400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -64300
$ws.Range("D96").Value = -70200
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -106800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -34700
Write-Output "done"
